$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.783.98'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.638.14'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '215.66'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").Value = '19.63'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  +1.63%  '
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '1.864.83'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '1.639.16'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '0.0₃0763'
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("D17").Value = '63.18'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '25.827.26'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '4.48'
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").Value = '192.49'
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = '9.96'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("E24").Value = '  +5.23%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = '141.68'
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("E27").Value = '  +1.37%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").Value = '15.51'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '0.0492'
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '0.905'
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").Value = '1.136.16'
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("E38").Value = '  -1.53%  '
$ws.Range("D39").Value = '0.546'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").Value = '100.62'
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").Value = '0.804'
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").Value = '1.773.52'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("E46").Value = '  +3.22%  '
$ws.Range("D47").Value = '55.31'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("E50").Value = '  +3.75%  '
$ws.Range("E51").Value = '  -2.44%  '
